{"js": "// Renumber [[PERSON_NN]] placeholder tags per the \"after fixing animal-name\n// surnames\" edit: one stray/duplicate PERSON id was removed from the\n// anonymization map, so every later id shifts down (most by 2, a handful by\n// 1) while two old ids (34 and 35) collapse onto the same new id (33).\n// Because the shift is position-dependent (not a uniform decrement), we\n// resolve it with an explicit old-id -> new-id table applied in a single\n// pass over each run of text (so an already-rewritten id is never\n// re-rewritten).\nconst personIdMap = {\n  \"27\": \"26\",\n  \"28\": \"27\",\n  \"29\": \"28\",\n  \"30\": \"29\",\n  \"31\": \"30\",\n  \"32\": \"31\",\n  \"33\": \"32\",\n  \"34\": \"33\",\n  \"35\": \"33\",\n  \"36\": \"34\",\n  \"37\": \"35\",\n  \"38\": \"36\",\n  \"39\": \"37\",\n  \"40\": \"38\",\n  \"41\": \"39\",\n  \"42\": \"40\",\n  \"43\": \"41\",\n  \"44\": \"42\",\n  \"45\": \"43\",\n  \"46\": \"44\",\n  \"47\": \"45\"\n};\n\nfunction remap(text) {\n  return text.replace(/\\[\\[PERSON_(\\d+)\\]\\]/g, function (match, num) {\n    return Object.prototype.hasOwnProperty.call(personIdMap, num)\n      ? \"[[PERSON_\" + personIdMap[num] + \"]]\"\n      : match;\n  });\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const oldText = paragraph.text;\n  if (oldText.indexOf(\"[[PERSON_\") === -1) {\n    continue;\n  }\n  const newText = remap(oldText);\n  if (newText !== oldText) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Renumber [[PERSON_NN]] placeholder tags per the \"after fixing animal-name\n# surnames\" edit: one stray/duplicate PERSON id was removed from the\n# anonymization map, so every later id shifts down (most by 2, a handful by\n# 1) while two old ids (34 and 35) collapse onto the same new id (33).\n# Because the shift is position-dependent (not a uniform decrement), we\n# resolve it with an explicit old-id -> new-id table applied in a single\n# pass over each paragraph's text (so an already-rewritten id is never\n# re-rewritten by a later rule).\n$personIdMap = @{\n    \"27\" = \"26\"\n    \"28\" = \"27\"\n    \"29\" = \"28\"\n    \"30\" = \"29\"\n    \"31\" = \"30\"\n    \"32\" = \"31\"\n    \"33\" = \"32\"\n    \"34\" = \"33\"\n    \"35\" = \"33\"\n    \"36\" = \"34\"\n    \"37\" = \"35\"\n    \"38\" = \"36\"\n    \"39\" = \"37\"\n    \"40\" = \"38\"\n    \"41\" = \"39\"\n    \"42\" = \"40\"\n    \"43\" = \"41\"\n    \"44\" = \"42\"\n    \"45\" = \"43\"\n    \"46\" = \"44\"\n    \"47\" = \"45\"\n}\n\nfunction Remap-PersonText($text) {\n    $re = [regex]'\\[\\[PERSON_(\\d+)\\]\\]'\n    $ms = $re.Matches($text)\n    if ($ms.Count -eq 0) {\n        return $text\n    }\n    $result = \"\"\n    $lastEnd = 0\n    foreach ($m in $ms) {\n        $result += $text.Substring($lastEnd, $m.Index - $lastEnd)\n        $num = $m.Groups[1].Value\n        if ($personIdMap.ContainsKey($num)) {\n            $result += \"[[PERSON_\" + $personIdMap[$num] + \"]]\"\n        } else {\n            $result += $m.Value\n        }\n        $lastEnd = $m.Index + $m.Length\n    }\n    $result += $text.Substring($lastEnd)\n    return $result\n}\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text\n    # Range.Text includes trailing paragraph / cell-end marks (chr 13 / chr 7);\n    # strip them before comparing/rewriting so we don't touch the marks.\n    $trimmed = $t.TrimEnd([char]13, [char]7)\n    if ($trimmed.Contains(\"[[PERSON_\")) {\n        $newText = Remap-PersonText $trimmed\n        if ($newText -ne $trimmed) {\n            $r.Text = $newText\n        }\n    }\n}\n\nWrite-Output \"done\"\n"}
